# Add a 'hole_id' index column (column A) to the "train" worksheet so that
# cross validation can match rows back to their original borehole IDs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("train")

$holeIds = @(
    "BRG_16_04A",
    "BRG_13_01",
    "BRG_16_03",
    "BRG_05_11",
    "ECO_09_03",
    "BRG_01_06",
    "BRG_16_04B",
    "ECO_09_04",
    "BRG_01_02",
    "BRG_05_13",
    "BRG_01_03",
    "BRG_05_12",
    "BRG_05_09",
    "BRG_01_08",
    "BRG_05_04",
    "BRG_05_15",
    "ECO_09_02",
    "BRG_01_07",
    "BRG_13_02",
    "ECO_09_01",
    "BRG_16_08",
    "BRG_05_01",
    "BRG_16_02",
    "BRG_05_03",
    "BRG_05_02",
    "BRG_05_14",
    "BRG_08_01",
    "BRG_01_01",
    "BRG_01_09",
    "BRG_01_04",
    "BRG_16_01"
)

# Header cell for the new index column - reuse the style already applied to
# the rest of column A / the header row (bold, centered, bordered) instead of
# building a brand-new style.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A1").Value = "hole_id"

for ($i = 0; $i -lt $holeIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $holeIds[$i]
}
